# Insert a new data row above row 228 (pushing the existing rows 228-304
# down to 229-305) and populate the new row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(228).Insert()

$ws.Range("A228").Value = 11
$ws.Range("B228").Value = "Vega Monumental Concepción"
$ws.Range("C228").Value = "Bíobío"
$ws.Range("D228").Value = 45120
$ws.Range("E228").Value = 8
$ws.Range("F228").Value = 100112003
$ws.Range("G228").Value = "Ajo"
$ws.Range("H228").Value = "Chino"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 200
$ws.Range("K228").Value = 16000
$ws.Range("L228").Value = 17000
$ws.Range("M228").Value = 16500
$ws.Range("N228").Value = "$/caja 10 kilos"
$ws.Range("O228").Value = "China"
$ws.Range("P228").Value = 1650
$ws.Range("Q228").Value = 10
$ws.Range("R228").Value = "Hortaliza"
